# Daily Report update: 2026-01-28
# Adds the new daily-data rows for date-serial 46049 to Daily_Data,
# and refreshes the dependent roll-up figures on Today_Summary and
# Monthly_Stats that shift because of JP MORGAN CHASE BANK NA's
# "Eligible" withdrawal recorded on the new date.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Daily_Data — append rows 354-375 (date serial 46049)
# ---------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    ,@(354, 'ASAHI DEPOSITORY LLC Registered', 0, 0, 0, 0, 0, 0)
    ,@(355, 'ASAHI DEPOSITORY LLC Eligible', 0, 0, 0, 0, 0, 0)
    ,@(356, 'BRINK''S, INC. Registered', 87949.747, 0, 0, 0, 0, 87949.747)
    ,@(357, 'BRINK''S, INC. Eligible', 30578.352, 0, 0, 0, 0, 30578.352)
    ,@(358, 'CNT DEPOSITORY, INC. Registered', 1246.06, 0, 0, 0, 0, 1246.06)
    ,@(359, 'CNT DEPOSITORY, INC. Eligible', 0, 0, 0, 0, 0, 0)
    ,@(360, 'DELAWARE DEPOSITORY Registered', 1633.941, 0, 0, 0, 0, 1633.941)
    ,@(361, 'DELAWARE DEPOSITORY Eligible', 18459.584, 0, 0, 0, 0, 18459.584)
    ,@(362, 'HSBC BANK, USA Registered', 1394.758, 0, 0, 0, 0, 1394.758)
    ,@(363, 'HSBC BANK, USA Eligible', 9281.978999999999, 0, 0, 0, 0, 9281.978999999999)
    ,@(364, 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 2395.448, 0, 0, 0, 0, 2395.448)
    ,@(365, 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 0, 0, 0, 0, 0, 0)
    ,@(366, 'JP MORGAN CHASE BANK NA Registered', 114985.579, 0, 0, 0, 0, 114985.579)
    ,@(367, 'JP MORGAN CHASE BANK NA Eligible', 135413.823, 0, 10006.15, -10006.15, 0, 125407.673)
    ,@(368, 'LOOMIS INTERNATIONAL (US) LLC Registered', 63745.991, 0, 0, 0, 0, 63745.991)
    ,@(369, 'LOOMIS INTERNATIONAL (US) LLC Eligible', 132077.206, 0, 0, 0, 0, 132077.206)
    ,@(370, 'MALCA-AMIT USA, LLC Registered', 395.145, 0, 0, 0, 0, 395.145)
    ,@(371, 'MALCA-AMIT USA, LLC Eligible', 0, 0, 0, 0, 0, 0)
    ,@(372, 'MANFRA, TORDELLA & BROOKES, LLC Registered', 50220.42, 0, 0, 0, 0, 50220.42)
    ,@(373, 'MANFRA, TORDELLA & BROOKES, LLC Eligible', 1271.373, 0, 0, 0, 0, 1271.373)
    ,@(374, 'STONEX PRECIOUS METALS LLC Registered', 14122.765, 0, 0, 0, 0, 14122.765)
    ,@(375, 'STONEX PRECIOUS METALS LLC Eligible', 16.075, 0, 0, 0, 0, 16.075)
)

$dateSerial = 46049

foreach ($row in $newRows) {
    $r = $row[0]

    $wsDaily.Cells.Item($r, 1).Value = $dateSerial
    # Match the date-stamp formatting used by the rest of column A.
    $wsDaily.Cells.Item($r, 1).NumberFormat = $wsDaily.Cells.Item($r - 1, 1).NumberFormat

    $wsDaily.Cells.Item($r, 2).Value = $row[1]
    $wsDaily.Cells.Item($r, 3).Value = $row[2]
    $wsDaily.Cells.Item($r, 4).Value = $row[3]
    $wsDaily.Cells.Item($r, 5).Value = $row[4]
    $wsDaily.Cells.Item($r, 6).Value = $row[5]
    $wsDaily.Cells.Item($r, 7).Value = $row[6]
    $wsDaily.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------
# 2) Today_Summary — JP MORGAN CHASE BANK NA eligible/total drop
#    by the 10006.15 withdrawn on the new date.
# ---------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("Today_Summary")
$wsToday.Range("B8").Value = 125407.673
$wsToday.Range("D8").Value = 240393.252

# ---------------------------------------------------------------
# 3) Monthly_Stats — monthly eligible total / grand total drop by
#    the same 10006.15, and the JP MORGAN CHASE BANK NA Eligible
#    monthly roll-up row now carries the withdrawn/total figures.
# ---------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")
$wsMonthly.Range("B2").Value = 317092.242
$wsMonthly.Range("D2").Value = 655182.0959999999

$wsMonthly.Range("D19").Value = 10006.15
$wsMonthly.Range("E19").Value = 125407.673
